# Auto-generated edit script applying cryptos list update (commit: Sun Oct 27 10:39:28 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers (e.g. "581.96") need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values (losing the original text-cell semantics used throughout column D).
$textFormatCells = @(
    "D5",
    "D6",
    "D8",
    "D14",
    "D18",
    "D19",
    "D20",
    "D23",
    "D30",
    "D31",
    "D35",
    "D36",
    "D44",
    "D45"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New cell values (address -> text), taken from the updated crypto price feed.
$updates = [ordered]@{
    "D2" = "67.095.88"
    "D3" = "2.468.30"
    "E3" = "  -0.08%  "
    "E4" = "  +0.00%  "
    "D5" = "581.96"
    "E5" = "  -0.24%  "
    "D6" = "173.95"
    "E6" = "  +2.32%  "
    "E7" = "  +0.02%  "
    "D8" = "0.512"
    "E8" = "  -0.49%  "
    "E9" = "  +1.49%  "
    "E10" = "  +0.18%  "
    "E11" = "  -0.66%  "
    "E12" = "  +0.48%  "
    "E13" = "  +0.03%  "
    "D14" = "25.34"
    "E14" = "  -1.01%  "
    "D15" = "66.965.17"
    "E16" = "  -0.13%  "
    "D17" = "2.420.92"
    "E17" = "  -1.79%  "
    "B18" = "Chainlink"
    "C18" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D18" = "10.87"
    "E18" = "  -2.38%  "
    "B19" = "Uniswap"
    "C19" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "D19" = "7.45"
    "E19" = "  -2.01%  "
    "D20" = "347.93"
    "E20" = "  -1.63%  "
    "E21" = "  -0.51%  "
    "D23" = "69.32"
    "E23" = "  +0.41%  "
    "E24" = "  -1.52%  "
    "E25" = "  +0.17%  "
    "E26" = "  -1.01%  "
    "D27" = "2.594.68"
    "E27" = "  +0.23%  "
    "E28" = "  +0.00%  "
    "D29" = "0.0₃0900"
    "E29" = "  -0.71%  "
    "D30" = "498.91"
    "E30" = "  -4.03%  "
    "D31" = "7.71"
    "E31" = "  -0.76%  "
    "E32" = "  -0.47%  "
    "E33" = "  -1.46%  "
    "E34" = "  -0.03%  "
    "D35" = "0.120"
    "E35" = "  +0.92%  "
    "D36" = "161.83"
    "E36" = "  +2.32%  "
    "E37" = "  -0.04%  "
    "E38" = "  -1.67%  "
    "E39" = "  -2.14%  "
    "E40" = "  -0.04%  "
    "E41" = "  +0.84%  "
    "E42" = "  -0.23%  "
    "E43" = "  -0.05%  "
    "D44" = "2.39"
    "E44" = "  +0.05%  "
    "D45" = "142.14"
    "E45" = "  +0.80%  "
    "E46" = "  +0.40%  "
    "E47" = "  -1.39%  "
    "E48" = "  -1.37%  "
    "E49" = "  +0.76%  "
    "E50" = "  -1.88%  "
    "E51" = "  -0.11%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Applied $($updates.Count) cell updates"
